{"js": "// Apply the dotNetRDF 0.5.0 design-document revision:\n//  - Target Date: November 2011 -> August 2011\n//  - \"Last Updated\" SAVEDATE field display text refreshed to 04/07/2011 11:07:00\n//  - _GoBack bookmark relocated from the title to the end of the date paragraph\n//  - Two new bullet points added under \"Required Features\"\n//  - New \"Additional Handlers\" section (heading + paragraph) added near the end\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------\n// 1) Target Date: November 2011 -> August 2011\n// ---------------------------------------------------------------\nconst dateSearch = body.search(\"November 2011\", { matchCase: true });\ndateSearch.load(\"items\");\nawait context.sync();\nif (dateSearch.items.length > 0) {\n  dateSearch.items[0].insertText(\"August 2011\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------\n// 2) Refresh the \"Last Updated\" SAVEDATE field display text\n// ---------------------------------------------------------------\nconst fieldSearch = body.search(\"14/04/2011 10:06:00\", { matchCase: true });\nfieldSearch.load(\"items\");\nawait context.sync();\nif (fieldSearch.items.length > 0) {\n  fieldSearch.items[0].insertText(\"04/07/2011 11:07:00\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------\n// 3) Move the _GoBack bookmark from the title paragraph to the end\n//    of the paragraph that holds the \"Last Updated\" date field\n// ---------------------------------------------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst allParas = body.paragraphs;\nallParas.load(\"items/text\");\nawait context.sync();\n\nlet dateParaIndex = -1;\nfor (let i = 0; i < allParas.items.length; i++) {\n  if (allParas.items[i].text.indexOf(\"Last Updated:\") !== -1) {\n    dateParaIndex = i;\n    break;\n  }\n}\nif (dateParaIndex !== -1) {\n  const endRange = allParas.items[dateParaIndex].getRange(Word.RangeLocation.end);\n  endRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// ---------------------------------------------------------------\n// 4) Add two new bullet points after \"Extend IRdfHandler ...\"\n//    under the \"Required Features\" heading\n// ---------------------------------------------------------------\nconst paras2 = body.paragraphs;\nparas2.load(\"items/text\");\nawait context.sync();\n\nlet extendIdx = -1;\nfor (let i = 0; i < paras2.items.length; i++) {\n  if (paras2.items[i].text.indexOf(\"Extend IRdfHandler and ISparqlResultsHandler\") !== -1) {\n    extendIdx = i;\n    break;\n  }\n}\nif (extendIdx !== -1) {\n  const anchor = paras2.items[extendIdx];\n  const newPara1 = anchor.insertParagraph(\n    \"Add additional handlers and associated required classes for serializing direct to other formats like RDF/XML, SPARQL XML etc\",\n    Word.InsertLocation.after\n  );\n  newPara1.styleBuiltIn = Word.Style.listParagraph;\n  await context.sync();\n  newPara1.attachToList(1, 0);\n  await context.sync();\n\n  const newPara2 = newPara1.insertParagraph(\n    \"Implement all features originally planned for the 0.4.2 and 0.4.3 releases, see relevant design documents for details\",\n    Word.InsertLocation.after\n  );\n  newPara2.styleBuiltIn = Word.Style.listParagraph;\n  await context.sync();\n  newPara2.attachToList(1, 0);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------\n// 5) Add the new \"Additional Handlers\" section using the first of\n//    the two trailing empty Heading 2 paragraphs\n// ---------------------------------------------------------------\nconst paras3 = body.paragraphs;\nparas3.load(\"items/text,items/style\");\nawait context.sync();\n\nlet emptyHeadingIdx = -1;\nfor (let i = 0; i < paras3.items.length; i++) {\n  if (paras3.items[i].style === \"Heading 2\" && paras3.items[i].text.trim() === \"\") {\n    emptyHeadingIdx = i;\n    break;\n  }\n}\nif (emptyHeadingIdx !== -1) {\n  const headingPara = paras3.items[emptyHeadingIdx];\n  headingPara.insertText(\"Additional Handlers\", Word.InsertLocation.end);\n  await context.sync();\n\n  const bodyPara = headingPara.insertParagraph(\n    \"Define new interfaces for allowing formatting of data into formats that require a header and/or footer such as XML based formats.  Alter existing write through handlers to support these.  Add a write through handler for SPARQL results\",\n    Word.InsertLocation.after\n  );\n  bodyPara.styleBuiltIn = Word.Style.normal;\n  await context.sync();\n}\n", "ps1": "# Apply the dotNetRDF 0.5.0 design-document revision:\n#  - Target Date: November 2011 -> August 2011\n#  - \"Last Updated\" SAVEDATE field display text refreshed to 04/07/2011 11:07:00\n#  - Two new bullet points added under \"Required Features\"\n#  - New \"Additional Handlers\" section (heading + paragraph) added near the end\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndex($doc, $substr) {\n    $paras = $doc.Paragraphs\n    for ($i = 1; $i -le $paras.Count; $i++) {\n        if ($paras.Item($i).Range.Text -like \"*$substr*\") {\n            return $i\n        }\n    }\n    return -1\n}\n\nfunction Find-EmptyHeading2Index($doc) {\n    $paras = $doc.Paragraphs\n    for ($i = 1; $i -le $paras.Count; $i++) {\n        $p = $paras.Item($i)\n        if ($p.Style.NameLocal -eq \"Heading 2\" -and $p.Range.Text.Trim() -eq \"\") {\n            return $i\n        }\n    }\n    return -1\n}\n\n# ---------------------------------------------------------------\n# 1) Target Date: November 2011 -> August 2011\n# ---------------------------------------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"November 2011\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"August 2011\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# ---------------------------------------------------------------\n# 2) Refresh the \"Last Updated\" SAVEDATE field display text\n# ---------------------------------------------------------------\nif ($d.Fields.Count -ge 1) {\n    $dateField = $d.Fields.Item(1)\n    $dateField.Result.Text = \"04/07/2011 11:07:00\"\n}\n\n# ---------------------------------------------------------------\n# 3) Add two new bullet points after \"Extend IRdfHandler ...\"\n#    under the \"Required Features\" heading. We duplicate the\n#    existing bullet paragraph (copy/paste) so the new paragraphs\n#    inherit the same list numbering (numId) and style.\n# ---------------------------------------------------------------\n$extendIdx = Find-ParagraphIndex $d \"Extend IRdfHandler and ISparqlResultsHandler\"\nif ($extendIdx -ne -1) {\n    $anchor = $d.Paragraphs.Item($extendIdx)\n    $anchor.Range.Copy()\n    $insertPoint = $d.Range($anchor.Range.End, $anchor.Range.End)\n    $insertPoint.Paste()\n\n    $newPara1 = $d.Paragraphs.Item($extendIdx + 1)\n    $r1 = $newPara1.Range\n    $r1.MoveEnd(1, -1) | Out-Null\n    $r1.Text = \"Add additional handlers and associated required classes for serializing direct to other formats like RDF/XML, SPARQL XML etc\"\n\n    $newPara1b = $d.Paragraphs.Item($extendIdx + 1)\n    $newPara1b.Range.Copy()\n    $insertPoint2 = $d.Range($newPara1b.Range.End, $newPara1b.Range.End)\n    $insertPoint2.Paste()\n\n    $newPara2 = $d.Paragraphs.Item($extendIdx + 2)\n    $r2 = $newPara2.Range\n    $r2.MoveEnd(1, -1) | Out-Null\n    $r2.Text = \"Implement all features originally planned for the 0.4.2 and 0.4.3 releases, see relevant design documents for details\"\n}\n\n# ---------------------------------------------------------------\n# 4) Add the new \"Additional Handlers\" section using the first of\n#    the two trailing empty Heading 2 paragraphs, followed by a\n#    new Normal-style paragraph (duplicated from an existing one\n#    so it carries the correct style).\n# ---------------------------------------------------------------\n$emptyHeadingIdx = Find-EmptyHeading2Index $d\nif ($emptyHeadingIdx -ne -1) {\n    $headingPara = $d.Paragraphs.Item($emptyHeadingIdx)\n    $rh = $headingPara.Range\n    $rh.MoveEnd(1, -1) | Out-Null\n    $rh.Text = \"Additional Handlers\"\n\n    $normalIdx = Find-ParagraphIndex $d \"Any method which executes a SPARQL Query\"\n    if ($normalIdx -ne -1) {\n        $normalPara = $d.Paragraphs.Item($normalIdx)\n        $normalPara.Range.Copy()\n\n        $headingPara2 = $d.Paragraphs.Item($emptyHeadingIdx)\n        $insertPoint3 = $d.Range($headingPara2.Range.End, $headingPara2.Range.End)\n        $insertPoint3.Paste()\n\n        $newBodyPara = $d.Paragraphs.Item($emptyHeadingIdx + 1)\n        $rb = $newBodyPara.Range\n        $rb.MoveEnd(1, -1) | Out-Null\n        $rb.Text = \"Define new interfaces for allowing formatting of data into formats that require a header and/or footer such as XML based formats.  Alter existing write through handlers to support these.  Add a write through handler for SPARQL results\"\n    }\n}\n"}
